$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph (right before
# the sectPr). We append seven new paragraphs after it:
#   1. (empty)
#   2. "TODO:"
#   3. "Homepage fertigstellen"
#   4. (empty)
#   5. "-Neue UI Page mit " / "Date" / "neingabe -- Stadtsuche:"  (3 runs)
#   6. "-> Verwendung von Tabelle"
#   7. "-> switch auf die Homepage zur Anzeige von Daten"
#
# Word's Range.InsertAfter() on a range collapsed to the very end of the
# document writes into the *current* last (trailing) paragraph - it does not
# create a fresh paragraph by itself.  So, for every new paragraph we first
# call InsertParagraphAfter() to open up a new trailing paragraph, and only
# then (for non-blank lines) fill it in with InsertAfter(text). Re-reading
# $d.Content fresh each time keeps us anchored at the true end of the
# document as it grows.

function Append-Paragraph {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

function Append-Text($text) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertAfter($text)
}

# 1. new empty paragraph
Append-Paragraph

# 2. "TODO:"
Append-Paragraph
Append-Text("TODO:")

# 3. "Homepage fertigstellen"
Append-Paragraph
Append-Text("Homepage fertigstellen")

# 4. empty paragraph
Append-Paragraph

# 5. three-run paragraph: "-Neue UI Page mit " + "Date" + "neingabe -- Stadtsuche:"
Append-Paragraph
Append-Text("-Neue UI Page mit ")
$afterFirstRun = $d.Content.End
Append-Text("Date")
$afterSecondRun = $d.Content.End
Append-Text("neingabe -- Stadtsuche:")

# 6. "-> Verwendung von Tabelle"
Append-Paragraph
Append-Text("-> Verwendung von Tabelle")

# 7. "-> switch auf die Homepage zur Anzeige von Daten"
Append-Paragraph
Append-Text("-> switch auf die Homepage zur Anzeige von Daten")

# Force the "Date" run (inserted in step 5) to stay a distinct <w:r> -
# matching the diff's 3-run split - by toggling a character property on
# just that sub-range *after* all the other text has already been typed.
# Toggling it back off again leaves the text/formatting unchanged but keeps
# the run boundary that plain adjoining InsertAfter calls would otherwise
# merge away; doing this last avoids leaking the touched-but-empty run
# properties into any later paragraphs.
$dateRange = $d.Range($afterFirstRun - 1, $afterSecondRun - 1)
$dateRange.Bold = 1
$dateRange.Bold = 0

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
